$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Wipe the old table contents (B1:J12) and drop the stray bold
#    formatting on the cells that become blank in the new layout.
# ------------------------------------------------------------------
$ws.Range("B1:J12").ClearContents()
$ws.Range("B2,B11,B12").Font.Bold = $false

# Old merged banner (E4:J4) goes away; a new one (B1:C1) appears.
$ws.Range("E4:J4").UnMerge()
$ws.Range("E4:J4").HorizontalAlignment = -4131

# ------------------------------------------------------------------
# 2) Title, merged across B1:C1, bold + centered.
# ------------------------------------------------------------------
$ws.Range("B1").Value = "Best r2 scores Copenhagen"
$ws.Range("B1:C1").Font.Bold = $true
$ws.Range("B1:C1").HorizontalAlignment = -4108
$ws.Range("B1:C1").Merge()

# Footnotes, now living at H1 / H2.
$ws.Range("H1").Value = "*LOOCV = leave-one-out cross validation"
$ws.Range("H2").Value = "*average = average score of model with 10 different train/test splits"

# ------------------------------------------------------------------
# 3) Column headers (row 3) + data rows, typed in roughly the same
#    order the original author entered them (keeps the shared-string
#    table close to the authored workbook).
# ------------------------------------------------------------------
$ws.Range("B3").Value = "Model"
$ws.Range("D3").Value = "Train r2"
$ws.Range("E3").Value = "Test r2"

$ws.Range("B4").Value = "Linear Regression"
$ws.Range("C4").Value = "LOOCV"

$ws.Range("B5").Value = "Linear Regression"
$ws.Range("C5").Value = "Average"
$ws.Range("D5").Value = 0.5155
$ws.Range("D5").Font.Name = "Calibri"
$ws.Range("E5").Value = 0.5147

$ws.Range("B6").Value = "NuSVR"
$ws.Range("C6").Value = "Average"
$ws.Range("D6").Value = 0.621
$ws.Range("E6").Value = 0.5363

$ws.Range("B8").Value = "Random Forest"
$ws.Range("C8").Value = "Average"

$ws.Range("B7").Value = "Gradient Boosting"
$ws.Range("C7").Value = "Average"

$ws.Range("B9").Value = "XGBoost"
$ws.Range("C9").Value = "Average"

$ws.Range("B10").Value = "Neural Network"
$ws.Range("C10").Value = "Average"

$ws.Range("C3").Value = "Eval method"
$ws.Range("B3:C3").Font.Bold = $true

# ------------------------------------------------------------------
# 5) "GBR 5000 datapoints" banner, merged across H4:M4.
# ------------------------------------------------------------------
$ws.Range("H4").Value = "GBR 5000 datapoints"
$ws.Range("H4:M4").Merge()

# ------------------------------------------------------------------
# 6) Column B width follows the new (shorter) labels - narrower than
#    before, auto-fit to the longest entry ("Linear Regression" /
#    "Gradient Boosting").
# ------------------------------------------------------------------
$ws.Columns("B").AutoFit()
$ws.Columns("B").ColumnWidth = 14.6

# ------------------------------------------------------------------
# 7) Selection cosmetics to match the saved view state.
# ------------------------------------------------------------------
$ws.Range("C4").Select()
